# Update countries & provincias Spain
# Applies the 6-Jul-2020 19:41 COVID-19 data refresh to the "Pais" sheet:
#  - refreshes the "last updated" timestamp
#  - updates Casos totales/Nuevos casos/Casos activos/Recuperados/Casos criticos/
#    Muertes hoy/Muertes for every row whose figures moved
#  - a handful of neighbouring countries swap rank (and so swap rows) once the
#    totals above are applied; their country-name cells are corrected too
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 6 de Julio de 2020 a las 19:41"

# Row 4
$ws.Range("B4").Value = 3005791
$ws.Range("C4").Value = 22863
$ws.Range("D4").Value = 1295042
$ws.Range("E4").Value = 1578065
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 115
$ws.Range("H4").Value = 132684

# Row 5
$ws.Range("B5").Value = 1613351
$ws.Range("C5").Value = 8766
$ws.Range("D5").Value = 978615
$ws.Range("E5").Value = 569616
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 220
$ws.Range("H5").Value = 65120

# Row 6
$ws.Range("B6").Value = 719401
$ws.Range("C6").Value = 21565
$ws.Range("D6").Value = 440099
$ws.Range("E6").Value = 259129
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 473
$ws.Range("H6").Value = 20173

# Row 17
$ws.Range("B17").Value = 206844
$ws.Range("C17").Value = 1086
$ws.Range("D17").Value = 182995
$ws.Range("E17").Value = 18608
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 5241

# Row 18
$ws.Range("B18").Value = 197888
$ws.Range("C18").Value = 330
$ws.Range("D18").Value = 182200
$ws.Range("E18").Value = 6602
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 9086

# Row 45
$ws.Range("B45").Value = 38128
$ws.Range("C45").Value = 703
$ws.Range("D45").Value = 19489
$ws.Range("E45").Value = 17835
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 10
$ws.Range("H45").Value = 804

# Row 49
$ws.Range("B49").Value = 30749
$ws.Range("C49").Value = 791
$ws.Range("D49").Value = 18056
$ws.Range("E49").Value = 12359
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 334

# Row 63
$ws.Range("A63").Value = "Argelia"
$ws.Range("B63").Value = 16404
$ws.Range("C63").Value = 463
$ws.Range("D63").Value = 11884
$ws.Range("E63").Value = 3561
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 7
$ws.Range("H63").Value = 959

# Row 64
$ws.Range("A64").Value = "Nepal"
$ws.Range("B64").Value = 15964
$ws.Range("C64").Value = 180
$ws.Range("D64").Value = 6811
$ws.Range("E64").Value = 9118
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 35

# Row 65
$ws.Range("B65").Value = 14379
$ws.Range("C65").Value = 164
$ws.Range("D65").Value = 10173
$ws.Range("E65").Value = 3969
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 237

# Row 79
$ws.Range("A79").Value = "Senegal"
$ws.Range("B79").Value = 7478
$ws.Range("C79").Value = 78
$ws.Range("D79").Value = 4870
$ws.Range("E79").Value = 2472
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = 136

# Row 80
$ws.Range("A80").Value = "Consejo Danes para los Refugiados"
$ws.Range("B80").Value = 7432
$ws.Range("C80").Value = 21
$ws.Range("D80").Value = 3226
$ws.Range("E80").Value = 4024
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 182

# Row 91
$ws.Range("B91").Value = 4996
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 1745
$ws.Range("E91").Value = 3229
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 3
$ws.Range("H91").Value = 22

# Row 99
$ws.Range("B99").Value = 3562
$ws.Range("C99").Value = 43
$ws.Range("D99").Value = 1374
$ws.Range("E99").Value = 1996
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 192

# Row 104
$ws.Range("B104").Value = 3048
$ws.Range("C104").Value = 51
$ws.Range("D104").Value = 1014
$ws.Range("E104").Value = 1942
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 92

# Row 126
$ws.Range("B126").Value = 1463
$ws.Range("C126").Value = 12
$ws.Range("D126").Value = 671
$ws.Range("E126").Value = 775
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 17

# Row 130
$ws.Range("B130").Value = 1199
$ws.Range("C130").Value = 11
$ws.Range("D130").Value = 1049
$ws.Range("E130").Value = 100
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 50

# Row 136
$ws.Range("A136").Value = "Mozambique"
$ws.Range("B136").Value = 1012
$ws.Range("C136").Value = 25
$ws.Range("D136").Value = 277
$ws.Range("E136").Value = 727
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 8

# Row 137
$ws.Range("A137").Value = "Republica de Chipre"
$ws.Range("B137").Value = 1003
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 839
$ws.Range("E137").Value = 145
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 19

# Row 138
$ws.Range("A138").Value = "Burkina Faso"
$ws.Range("B138").Value = 1000
$ws.Range("C138").Value = 13
$ws.Range("D138").Value = 858
$ws.Range("E138").Value = 89
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 53

# Row 139
$ws.Range("A139").Value = "Suazilandia"
$ws.Range("B139").Value = 988
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 547
$ws.Range("E139").Value = 428
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 13

# Row 146
$ws.Range("B146").Value = 841
$ws.Range("C146").Value = 60
$ws.Range("D146").Value = 315
$ws.Range("E146").Value = 512
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 14

# Row 165
$ws.Range("A165").Value = "Botsuana"
$ws.Range("B165").Value = 314
$ws.Range("C165").Value = 37
$ws.Range("D165").Value = 31
$ws.Range("E165").Value = 282
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 1

# Row 166
$ws.Range("A166").Value = "Comoras"
$ws.Range("B166").Value = 311
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 266
$ws.Range("E166").Value = 38
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 7

# Row 184
$ws.Range("A184").Value = "Lesoto"
$ws.Range("B184").Value = 91
$ws.Range("C184").Value = 12
$ws.Range("D184").Value = 11
$ws.Range("E184").Value = 80
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 0

# Row 185
$ws.Range("A185").Value = "Liechtenstein"
$ws.Range("B185").Value = 84
$ws.Range("C185").Value = 1
$ws.Range("D185").Value = 81
$ws.Range("E185").Value = 2
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 1

# Row 186
$ws.Range("A186").Value = "Seychelles"
$ws.Range("B186").Value = 81
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 11
$ws.Range("E186").Value = 70
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0

# Row 187
$ws.Range("A187").Value = "Butan"
$ws.Range("B187").Value = 80
$ws.Range("C187").Value = 2
$ws.Range("D187").Value = 53
$ws.Range("E187").Value = 27
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0

# Rows 209-210: Groenlandia and Islas Malvinas are tied on every figure, so only
# their country names trade places (no B:H values change)
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

